# Add a new knowledge-base article row (row 33) to Sheet1:
#   "Sensus Navigation (2016) 업데이트 하기" / HT203026 / <long HTML howto> / 2024-06-11 (45454)
# plus minor row-height touch-ups that a newer Excel build recomputed when the
# file was re-saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Minor row-height jitter (auto layout recompute from a newer Excel build).
# Purely cosmetic - does not change any cell content.
# ---------------------------------------------------------------------------
$ws.Rows(2).RowHeight = 95.1
$ws.Rows(3).RowHeight = 75.95
$ws.Rows(5).RowHeight = 75.95
$ws.Rows(6).RowHeight = 95.1
$ws.Rows(7).RowHeight = 75.95
$ws.Rows(8).RowHeight = 189.95
$ws.Rows(13).RowHeight = 75.95
$ws.Rows(14).RowHeight = 75.95
$ws.Rows(16).RowHeight = 95.1
$ws.Rows(17).RowHeight = 95.1
$ws.Rows(19).RowHeight = 95.1
$ws.Rows(20).RowHeight = 95.1
$ws.Rows(22).RowHeight = 323.10000000000002
$ws.Rows(23).RowHeight = 266.10000000000002
$ws.Rows(26).RowHeight = 152.1
$ws.Rows(27).RowHeight = 132.94999999999999

# ---------------------------------------------------------------------------
# New row 33: Sensus Navigation (2016) update how-to article
# ---------------------------------------------------------------------------
$title = 'Sensus Navigation (2016) 업데이트 하기'
$code = 'HT203026'
$desc = @'
Sensus Navigation(2016)이 구동 중인 경우 업데이트 방법에 대해 고객이 보고할 수 있습니다.
<br>
<br>차량이 네트워크가 동작하는 경우, 무선(OTA)로 업데이트할 수 있습니다.
<br><h3>인터넷을 통해 업데이트 방법
<br>1. 소프트웨어 업데이트는 운전 중이거나 차량의 시동이 켜져 있는 경우 차량에 직접 다운로드됩니다.
<br>2. 업데이트가 다운로드되면 업데이트를 설치할 준비가 되었다는 알림을 받게 됩니다. '지금 설치'를 클릭하여 설치를 시작합니다.
<br>즉시 설치하거나 원하는 시간을 선택하여 설치를 시작할 수 있습니다. 설치를 완료하는 데 약 90분이 소요됩니다.
<br>3.소프트웨가 최신 상태입니다. 메시지가 중앙 화면에 나타나면 소프트웨어 업데이트가 완료된 것입니다.
<br><h4>차량에서 직접 다운로드하기</h4>
<br>1. 앱 보기에서 다운로드 센터 버튼을 누릅니다.
<br>2. 지도 버튼을 누릅니다.
<br>3.설치 버튼을 누른 후 확인을 선택합니다.
<br>4. 선택한 지도 업데이트의 설치가 시작됩니다.
<br>
<br><h4>진행 중인 지도 다운로드가 취소될 경우
<br>지도 다운로드 중에 차량 시동을 끄면, 차량을 다시 시동하여 인터넷에 다시 연결할 때 지도 다운로드가 다시 시작됩니다.
<br><h3>지도를 USB로 설치하는 방법
<br>새 지도는 인터넷에 연결된 컴퓨터에서 USB 메모리로 다운로드한 후 USB 메모리에서 차량의 내비게이션 시스템으로 전송할 수 있습니다.
<br>
<br>USB 메모리
<br>USB 메모리로 업데이트를 관리해야 하는 경우에 다음 요구 사항이 적용됩니다.
<br>USB 표준:	최소 2.0
<br>파일 시스템: FAT32, exFAT 또는 NFTS
<br>용량: 최대 128GB
<br>
<br>지도 파일은 아래 링크에서 다운로드 받을 수 있습니다.
<br><a href="https://www.volvocars.com/kr/support/downloads/maps/spa/daehanmingug" target="_blank">Sensus Navigation 2016 지도 다운로드</a>
'@

$ws.Range("A33").Value = $title
$ws.Range("B33").Value = 45454
$ws.Range("C33").Value = $code
$ws.Range("D33").Value = $desc

# D column in this sheet uses a wrapped, vertically-centred style (same as
# the other "Data" cells, e.g. D32) - match it so D33 shares that style.
$ws.Range("D33").WrapText = $true
$ws.Range("D33").VerticalAlignment = -4108

$ws.Rows(33).RowHeight = 409.6

# Leave the new row's Category (E) column empty, matching the source data.

# Match the final active selection/view reflected in the saved workbook.
$ws.Range("D33").Select() | Out-Null
